{"js": "// The document has a short paragraph containing the text \"Teste\" followed\n// by a non-breaking space, with a \"_GoBack\" bookmark sitting between the\n// two runs. The edit removes both runs of text (the word \"Teste\" and the\n// trailing space) while leaving the bookmark (and the now-empty paragraph)\n// in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose entire (trimmed) text is exactly \"Teste\" --\n// this is the lone paragraph that holds the bookmark between the two runs\n// being removed (other occurrences of \"Teste\" in the document are part of\n// longer sentences, so a trimmed full-text match singles out the right one).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  // Trim regular spaces as well as non-breaking spaces (\\u00A0).\n  const text = para.text.replace(/[\\s\\u00A0]+$/g, \"\").replace(/^[\\s\\u00A0]+/g, \"\");\n  if (text === \"Teste\") {\n    target = para;\n    break;\n  }\n}\n\nif (target) {\n  // Remove the \"Teste\" run.\n  const wordRanges = target.search(\"Teste\", { matchCase: true });\n  wordRanges.load(\"items\");\n  await context.sync();\n  if (wordRanges.items.length > 0) {\n    wordRanges.items[0].delete();\n    await context.sync();\n  }\n\n  // Remove the remaining trailing non-breaking space run, leaving the\n  // bookmark (_GoBack) untouched inside the now text-less paragraph.\n  const spaceRanges = target.search(\"\\u00A0\", { matchCase: true });\n  spaceRanges.load(\"items\");\n  await context.sync();\n  if (spaceRanges.items.length > 0) {\n    spaceRanges.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The document has a short paragraph whose entire text is the word \"Teste\"\n# followed by a non-breaking space, with a \"_GoBack\" bookmark sitting\n# between the two text runs. The edit removes both pieces of text (the\n# word \"Teste\" and the trailing non-breaking space) while leaving the\n# bookmark - and the now text-less paragraph - in place.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose trimmed text is exactly \"Teste\". Other\n# occurrences of \"Teste\" in the document are part of longer sentences, so\n# matching the full (trimmed) paragraph text singles out the right one.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $t2 = $t.TrimEnd([char]13, [char]10, [char]32, [char]160)\n    $t2 = $t2.TrimStart([char]32, [char]160)\n    if ($t2 -eq \"Teste\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $pr = $target.Range\n\n    # Delete the \"Teste\" run (the first 5 characters of the paragraph).\n    $wordRange = $d.Range($pr.Start, $pr.Start + 5)\n    $wordRange.Text = \"\"\n\n    # Re-fetch the paragraph range and delete whatever text remains before\n    # the paragraph mark (the trailing non-breaking space run), leaving the\n    # bookmark untouched.\n    $pr2 = $target.Range\n    if (($pr2.End - 1) -gt $pr2.Start) {\n        $restRange = $d.Range($pr2.Start, $pr2.End - 1)\n        $restRange.Text = \"\"\n    }\n}\n"}
